$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily ridership rows (Tue 06 Jun 2017 .. Sun 11 Jun 2017)
$rows = @(
    @("Tuesday",   "06 Jun 2017", 234, 219.09, 123.53),
    @("Wednesday", "07 Jun 2017", 241, 215.84, 123.7),
    @("Thursday",  "08 Jun 2017", 240, 234.97, 123.86),
    @("Friday",    "09 Jun 2017", 257, 238.47, 124.02),
    @("Saturday",  "10 Jun 2017", 133, 114.82, 124.18),
    @("Sunday",    "11 Jun 2017", 94,  90.42,  124.34)
)

$r = 3
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

$lastRow = $r - 1

# Update the three chart series (Ridership, Average, Pilot) to cover the
# full new data range instead of just the single original row.
$co = $ws.ChartObjects(1)
$chart = $co.Chart

$ser1 = $chart.SeriesCollection(1)
$ser1.Formula = '=SERIES("Ridership",Ridership!$B$2:$B$' + $lastRow + ',Ridership!$C$2:$C$' + $lastRow + ',1)'

$ser2 = $chart.SeriesCollection(2)
$ser2.Formula = '=SERIES("Average",Ridership!$B$2:$B$' + $lastRow + ',Ridership!$D$2:$D$' + $lastRow + ',2)'

$ser3 = $chart.SeriesCollection(3)
$ser3.Formula = '=SERIES("Pilot",Ridership!$B$2:$B$' + $lastRow + ',Ridership!$E$2:$E$' + $lastRow + ',3)'

# The chart's anchor shifts down by the same 6 rows that were inserted
# (row 4 -> 10, row 18 -> 24), keeping the same size/offsets.
$co.Top = $co.Top + 6 * $ws.Rows.Item(1).RowHeight
